$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.112114191055298
$ws.Range("B1").Value = 2.998554468154907
$ws.Range("C1").Value = 2.545324087142944
$ws.Range("D1").Value = 2.461684465408325
$ws.Range("E1").Value = 1.932547450065613
